# "method for regular proposal started"
# - Rename the two "proposal" labels on Sheet1 (A2/B2) to their new NEW-SPS / NEW SPS values
# - Move the active selection on Sheet1 from B2 to C3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "NEW-SPS"
$ws.Range("B2").Value = "NEW SPS"

$ws.Activate()
$ws.Range("C3").Select()
